$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts rows 9-100 down to 10-101),
# creating room for a new daily entry (day 8 of 06/2025) that was
# added ahead of the existing 05/2025 block.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data point.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 4942.9
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 2025
$ws.Range("E9").Value = "06/2025"
